$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 23,10
$arr[0,0] = -0.8260246319213993
$arr[0,1] = -1.981077477814098
$arr[0,2] = -0.4414436339245075
$arr[0,3] = -0.7110977538981412
$arr[0,4] = -0.159566049999028
$arr[0,5] = -0.3854928109118805
$arr[0,6] = -0.2776770955153309
$arr[0,7] = -0.3740767161796326
$arr[0,8] = 0.1293341692733759
$arr[0,9] = 0.3878561111968251
$arr[1,0] = -0.1405287498260583
$arr[1,1] = -0.4101828697996921
$arr[1,2] = 0.1413488340994211
$arr[1,3] = -0.08457792681343129
$arr[1,4] = 0.0232377885831183
$arr[1,5] = -0.0731618320811834
$arr[1,6] = 0.4302490533718251
$arr[1,7] = 0.6887709952952742
$arr[1,8] = -1.03521158968304
$arr[1,9] = -0.3815723733203383
$arr[2,0] = 0.7815531178611421
$arr[2,1] = 0.5556263569482897
$arr[2,2] = 0.6634420723448393
$arr[2,3] = 0.5670424516805376
$arr[2,4] = 1.070453337133546
$arr[2,5] = 1.328975279056995
$arr[2,6] = -0.3950073059213186
$arr[2,7] = 0.2586319104413827
$arr[2,8] = 0.6221758453498138
$arr[2,9] = -0.2785999995298156
$arr[3,0] = 0.59085446968239
$arr[3,1] = 0.4944548490180882
$arr[3,2] = 0.9978657344710967
$arr[3,3] = 1.256387676394546
$arr[3,4] = -0.467594908583768
$arr[3,5] = 0.1860443077789333
$arr[3,6] = 0.5495882426873645
$arr[3,7] = -0.351187602192265
$arr[3,8] = 0.7191879547613879
$arr[3,9] = 0.142128040744926
$arr[4,0] = 0.9421746901075616
$arr[4,1] = 1.200696632031011
$arr[4,2] = -0.5232859529473031
$arr[4,3] = 0.1303532634153982
$arr[4,4] = 0.4938971983238294
$arr[4,5] = -0.4068786465558001
$arr[4,6] = 0.6634969103978527
$arr[4,7] = 0.08643699638139091
$arr[4,8] = -0.1191689877116253
$arr[4,9] = 0.8599759940094719
$arr[5,0] = -0.8773451131212686
$arr[5,1] = -0.2237058967585673
$arr[5,2] = 0.1398380381498639
$arr[5,3] = -0.7609378067297656
$arr[5,4] = 0.3094377502238873
$arr[5,5] = -0.2676221637925746
$arr[5,6] = -0.4732281478855908
$arr[5,7] = 0.5059168338355065
$arr[5,8] = -0.2304295095018555
$arr[5,9] = -0.2483993241560489
$arr[6,0] = 0.04951982135029001
$arr[6,1] = -0.8512560235293395
$arr[6,2] = 0.2191195334243133
$arr[6,3] = -0.3579403805921485
$arr[6,4] = -0.5635463646851647
$arr[6,5] = 0.4155986170359325
$arr[6,6] = -0.3207477263014294
$arr[6,7] = -0.3387175409556228
$arr[6,8] = -0.315312675415286
$arr[6,9] = -0.6804871111820268
$arr[7,0] = 0.2874422675287736
$arr[7,1] = -0.2896176464876882
$arr[7,2] = -0.4952236305807045
$arr[7,3] = 0.4839213511403928
$arr[7,4] = -0.2524249921969692
$arr[7,5] = -0.2703948068511625
$arr[7,6] = -0.2469899413108257
$arr[7,7] = -0.6121643770775664
$arr[7,8] = 0.04365514009338567
$arr[7,9] = -0.08721574084803801
$arr[8,0] = -0.5938852086106097
$arr[8,1] = 0.3852597731104875
$arr[8,2] = -0.3510865702268745
$arr[8,3] = -0.3690563848810678
$arr[8,4] = -0.345651519340731
$arr[8,5] = -0.7108259551074718
$arr[8,6] = -0.0550064379365196
$arr[8,7] = -0.1858773188779433
$arr[8,8] = -0.2953951178020652
$arr[8,9] = -1.125172904869618
$arr[9,0] = -0.2740769919852016
$arr[9,1] = -0.292046806639395
$arr[9,2] = -0.2686419410990581
$arr[9,3] = -0.6338163768657989
$arr[9,4] = 0.02200314030515318
$arr[9,5] = -0.1088677406362705
$arr[9,6] = -0.2183855395603924
$arr[9,7] = -1.048163326627946
$arr[9,8] = -0.4676375374592979
$arr[9,9] = -0.5510236388890779
$arr[10,0] = -0.09712588908289738
$arr[10,1] = -0.4623003248496382
$arr[10,2] = 0.193519192321314
$arr[10,3] = 0.0626483113798903
$arr[10,4] = -0.04686948754423159
$arr[10,5] = -0.8766472746117848
$arr[10,6] = -0.2961214854431372
$arr[10,7] = -0.3795075868729171
$arr[10,8] = -0.5974174644126757
$arr[10,9] = 0.9663603392234945
$arr[11,0] = 0.1560540901775642
$arr[11,1] = 0.02518320923614054
$arr[11,2] = -0.08433458968798135
$arr[11,3] = -0.9141123767555346
$arr[11,4] = -0.3335865875868869
$arr[11,5] = -0.4169726890166669
$arr[11,6] = -0.6348825665564255
$arr[11,7] = 0.9288952370797448
$arr[11,8] = -0.4070884690023188
$arr[11,9] = 1.033525947527809
$arr[12,0] = -0.106237102096024
$arr[12,1] = -0.9360148891635773
$arr[12,2] = -0.3554890999949296
$arr[12,3] = -0.4388752014247095
$arr[12,4] = -0.6567850789644682
$arr[12,5] = 0.9069927246717022
$arr[12,6] = -0.4289909814103615
$arr[12,7] = 1.011623435119767
$arr[12,8] = -0.0888206791410559
$arr[12,9] = -0.2990950773107258
$arr[13,0] = -0.2673832081489602
$arr[13,1] = -0.3507693095787401
$arr[13,2] = -0.5686791871184989
$arr[13,3] = 0.9950986165176716
$arr[13,4] = -0.3408850895643921
$arr[13,5] = 1.099729326965736
$arr[13,6] = -0.0007147872950865053
$arr[13,7] = -0.2109891854647564
$arr[13,8] = 0.387555017918974
$arr[13,9] = 0.06188637714375961
$arr[14,0] = -0.4749660705519546
$arr[14,1] = 1.088811733084216
$arr[14,2] = -0.2471719729978479
$arr[14,3] = 1.19344244353228
$arr[14,4] = 0.09299832927145768
$arr[14,5] = -0.1172760688982122
$arr[14,6] = 0.4812681344855182
$arr[14,7] = 0.1555994937103038
$arr[14,8] = 0.8727792950859075
$arr[14,9] = 2.568651570723608
$arr[15,0] = 1.097015279354651
$arr[15,1] = -0.2389684267274125
$arr[15,2] = 1.201645989802716
$arr[15,3] = 0.1012018755418931
$arr[15,4] = -0.1090725226277768
$arr[15,5] = 0.4894716807559536
$arr[15,6] = 0.1638030399807392
$arr[15,7] = 0.8809828413563429
$arr[15,8] = 2.576855116994044
$arr[15,9] = 9.421101911918621
$arr[16,0] = -0.2415579873788807
$arr[16,1] = 1.199056429151248
$arr[16,2] = 0.09861231489042488
$arr[16,3] = -0.111662083279245
$arr[16,4] = 0.4868821201044854
$arr[16,5] = 0.161213479329271
$arr[16,6] = 0.8783932807048747
$arr[16,7] = 2.574265556342576
$arr[16,8] = 9.418512351267154
$arr[16,9] = -8.132141731834157
$arr[17,0] = 1.205589761734299
$arr[17,1] = 0.1051456474734768
$arr[17,2] = -0.1051287506961931
$arr[17,3] = 0.4934154526875373
$arr[17,4] = 0.1677468119123229
$arr[17,5] = 0.8849266132879265
$arr[17,6] = 2.580798888925627
$arr[17,7] = 9.425045683850206
$arr[17,8] = -8.125608399251105
$arr[17,9] = -0.4162117995949584
$arr[18,0] = -0.006071047505593896
$arr[18,1] = -0.2163454456752638
$arr[18,2] = 0.3821987577084666
$arr[18,3] = 0.05653011693325222
$arr[18,4] = 0.7737099183088558
$arr[18,5] = 2.469582193946557
$arr[18,6] = 9.313828988871135
$arr[18,7] = -8.236825094230175
$arr[18,8] = -0.5274284945740291
$arr[18,9] = 1.055484166312883
$arr[19,0] = -0.2580915896621678
$arr[19,1] = 0.3404526137215625
$arr[19,2] = 0.01478397294634815
$arr[19,3] = 0.7319637743219518
$arr[19,4] = 2.427836049959653
$arr[19,5] = 9.272082844884231
$arr[19,6] = -8.278571238217079
$arr[19,7] = -0.5691746385609331
$arr[19,8] = 1.013738022325978
$arr[19,9] = -1.952548970023277
$arr[20,0] = 0.3657676764542774
$arr[20,1] = 0.04009903567906303
$arr[20,2] = 0.7572788370546667
$arr[20,3] = 2.453151112692368
$arr[20,4] = 9.297397907616945
$arr[20,5] = -8.253256175484365
$arr[20,6] = -0.5438595758282182
$arr[20,7] = 1.039053085058693
$arr[20,8] = -1.927233907290562
$arr[20,9] = 0.3367972473739005
$arr[21,0] = 0.04413770072197692
$arr[21,1] = 0.7613175020975806
$arr[21,2] = 2.457189777735282
$arr[21,3] = 9.30143657265986
$arr[21,4] = -8.249217510441451
$arr[21,5] = -0.5398209107853043
$arr[21,6] = 1.043091750101607
$arr[21,7] = -1.923195242247649
$arr[21,8] = 0.3408359124168144
$arr[21,9] = -0.1103182425099242
$arr[22,0] = 0.6424403654065582
$arr[22,1] = 2.338312641044259
$arr[22,2] = 9.182559435968837
$arr[22,3] = -8.368094647132473
$arr[22,4] = -0.6586980474763267
$arr[22,5] = 0.9242146134105849
$arr[22,6] = -2.042072378938671
$arr[22,7] = 0.221958775725792
$arr[22,8] = -0.2291953792009466
$arr[22,9] = -0.3695594427149207

$ws.Range("B2:K24").Value = $arr
